$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Jose Alvarado"
$ws.Range("B2").Value = "PG"
$ws.Range("C2").Value = "New Orleans Pelicans"

$ws.Range("A3").Value = "Shai Gilgeous-Alexander"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "Oklahoma City Thunder"

$ws.Range("A4").Value = "Luke Kennard"
$ws.Range("B4").Value = "SG"
$ws.Range("C4").Value = "Memphis Grizzlies"

$ws.Range("A5").Value = "CJ McCollum"
$ws.Range("B5").Value = "PG,SG"
$ws.Range("C5").Value = "New Orleans Pelicans"

$ws.Range("A6").Value = "Tobias Harris"
$ws.Range("B6").Value = "SF,PF"
$ws.Range("C6").Value = "Detroit Pistons"

$ws.Range("A7").Value = "Toumani Camara"
$ws.Range("B7").Value = "SF,PF"
$ws.Range("C7").Value = "Portland Trail Blazers"

$ws.Range("A8").Value = "Lauri Markkanen"
$ws.Range("B8").Value = "SF,PF"
$ws.Range("C8").Value = "Utah Jazz"

$ws.Range("A9").Value = "Jordan Clarkson"
$ws.Range("B9").Value = "SG,SF"
$ws.Range("C9").Value = "Utah Jazz"

$ws.Range("A10").Value = "John Collins"
$ws.Range("B10").Value = "PF,C"
$ws.Range("C10").Value = "Utah Jazz"

$ws.Range("A11").Value = "Wendell Carter Jr."
$ws.Range("B11").Value = "PF,C"
$ws.Range("C11").Value = "Orlando Magic"

$ws.Range("A12").Value = "Jordan Poole"
$ws.Range("B12").Value = "PG,SG"
$ws.Range("C12").Value = "Washington Wizards"

$ws.Range("A13").Value = "Jamal Murray"
$ws.Range("B13").Value = "PG,SG"
$ws.Range("C13").Value = "Denver Nuggets"

$ws.Range("A14").Value = "Quentin Grimes"
$ws.Range("B14").Value = "SG,SF"
$ws.Range("C14").Value = "Dallas Mavericks"

$ws.Range("A15").Value = "Bam Adebayo"
$ws.Range("B15").Value = "PF,C"
$ws.Range("C15").Value = "Miami Heat"

$ws.Range("A16").Value = "Nikola Jovic"
$ws.Range("B16").Value = "PF,C"
$ws.Range("C16").Value = "Miami Heat"

$ws.Range("A17").Value = "Zach LaVine"
$ws.Range("B17").Value = "SG,SF"
$ws.Range("C17").Value = "Sacramento Kings"

$ws.Range("A18").Value = "Joel Embiid"
$ws.Range("B18").Value = "C"
$ws.Range("C18").Value = "Philadelphia 76ers"

$ws.Range("A19").Value = "Kyrie Irving"
$ws.Range("B19").Value = "PG,SG"
$ws.Range("C19").Value = "Dallas Mavericks"
